# Applies the scheduled-runner market-data refresh described in the commit.
# For each affected Leve row, currentAveragePrice/NQ/HQ, LevePriceNQ/HQ and the
# derived LeveProfitNQ/HQ columns (H, I, J, K, L, M, N) are updated to the latest
# pulled values. Where the new computed profit is blank/not-applicable the cell
# is cleared entirely (matches the workbook's sparse convention of omitting 0/blank).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32
$ws.Range("H32").Value = 966.3333
$ws.Range("J32").Value = 966.3333
$ws.Range("L32").Value = 966.3333
$ws.Range("N32").Value = -1618.3333
# Row 51
$ws.Range("H51").Value = 1500
$ws.Range("I51").Value = 0
$ws.Range("K51").Value = 0
$ws.Range("M51").Value = ""
# Row 64
$ws.Range("H64").Value = 16500.25
$ws.Range("J64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("N64").Value = -5496
# Row 67
$ws.Range("H67").Value = 16500.25
$ws.Range("J67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("N67").Value = -6716
# Row 69
$ws.Range("H69").Value = 169766.33
$ws.Range("I69").Value = 2300
$ws.Range("J69").Value = 253499.5
$ws.Range("K69").Value = 6900
$ws.Range("L69").Value = 760498.5
$ws.Range("M69").Value = -6026
$ws.Range("N69").Value = -762246.5
# Row 72
$ws.Range("H72").Value = 169766.33
$ws.Range("I72").Value = 2300
$ws.Range("J72").Value = 253499.5
$ws.Range("K72").Value = 20700
$ws.Range("L72").Value = 2281495.5
$ws.Range("M72").Value = -16332
$ws.Range("N72").Value = -2290231.5
# Row 88
$ws.Range("H88").Value = 6655.875
$ws.Range("I88").Value = 3600
$ws.Range("J88").Value = 8489.4
$ws.Range("K88").Value = 3600
$ws.Range("L88").Value = 8489.4
$ws.Range("M88").Value = -3194
$ws.Range("N88").Value = -9301.4
# Row 91
$ws.Range("H91").Value = 6655.875
$ws.Range("I91").Value = 3600
$ws.Range("J91").Value = 8489.4
$ws.Range("K91").Value = 3600
$ws.Range("L91").Value = 8489.4
$ws.Range("M91").Value = -2196
$ws.Range("N91").Value = -11297.4
# Row 92
$ws.Range("H92").Value = 487
$ws.Range("I92").Value = 487
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 487
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = 761
$ws.Range("N92").Value = ""
# Row 99
$ws.Range("H99").Value = 200
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 200
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = ""
$ws.Range("N99").Value = -3596

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4096.2
$ws.Range("I32").Value = 4096.2
$ws.Range("K32").Value = 4096.2
$ws.Range("M32").Value = -3809.2
# Row 97
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = ""
$ws.Range("N97").Value = ""

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 5899.7646
$ws.Range("I86").Value = 3142.7144
$ws.Range("J86").Value = 7829.7
$ws.Range("K86").Value = 3142.7144
$ws.Range("L86").Value = 7829.7
$ws.Range("M86").Value = -2019.7144
$ws.Range("N86").Value = -10075.7
# Row 89
$ws.Range("H89").Value = 5899.7646
$ws.Range("I89").Value = 3142.7144
$ws.Range("J89").Value = 7829.7
$ws.Range("K89").Value = 15713.572
$ws.Range("L89").Value = 39148.5
$ws.Range("M89").Value = -10097.572
$ws.Range("N89").Value = -50380.5
# Row 94
$ws.Range("H94").Value = 6351.5
$ws.Range("I94").Value = 4552.8335
$ws.Range("J94").Value = 9049.5
$ws.Range("K94").Value = 4552.8335
$ws.Range("L94").Value = 9049.5
$ws.Range("M94").Value = -4101.8335
$ws.Range("N94").Value = -9951.5
# Row 99
$ws.Range("H99").Value = 3371.8
$ws.Range("I99").Value = 3214.75
$ws.Range("K99").Value = 3214.75
$ws.Range("M99").Value = -1716.75
# Row 134
$ws.Range("H134").Value = 3749
$ws.Range("I134").Value = 3436.25
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 10308.75
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -7773.75
$ws.Range("N134").Value = -20070

$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 161
$ws.Range("I22").Value = 252.5
$ws.Range("K22").Value = 252.5
$ws.Range("M22").Value = 97.5
# Row 31
$ws.Range("H31").Value = 2877.1765
$ws.Range("I31").Value = 1704.1111
$ws.Range("J31").Value = 4196.875
$ws.Range("K31").Value = 1704.1111
$ws.Range("L31").Value = 4196.875
$ws.Range("M31").Value = -1409.1111
$ws.Range("N31").Value = -4786.875
# Row 34
$ws.Range("H34").Value = 2877.1765
$ws.Range("I34").Value = 1704.1111
$ws.Range("J34").Value = 4196.875
$ws.Range("K34").Value = 1704.1111
$ws.Range("L34").Value = 4196.875
$ws.Range("M34").Value = -1502.1111
$ws.Range("N34").Value = -4600.875
# Row 132
$ws.Range("H132").Value = 104109.6
$ws.Range("I132").Value = 128637.25
$ws.Range("K132").Value = 385911.75
$ws.Range("M132").Value = -383381.75

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 992.5
$ws.Range("I5").Value = 788.44446
$ws.Range("J5").Value = 1359.8
$ws.Range("K5").Value = 2365.33338
$ws.Range("L5").Value = 4079.4
$ws.Range("M5").Value = -2253.33338
$ws.Range("N5").Value = -4303.4
# Row 34
$ws.Range("H34").Value = 2113.5
$ws.Range("J34").Value = 3825
$ws.Range("L34").Value = 11475
$ws.Range("N34").Value = -11643
# Row 37
$ws.Range("H37").Value = 150000
$ws.Range("J37").Value = 150000
$ws.Range("L37").Value = 450000
$ws.Range("N37").Value = -450224
# Row 92
$ws.Range("H92").Value = 350
$ws.Range("J92").Value = 350
$ws.Range("L92").Value = 1050
$ws.Range("N92").Value = -3546
# Row 99
$ws.Range("H99").Value = 4744.44
$ws.Range("J99").Value = 4822.174
$ws.Range("L99").Value = 14466.522
$ws.Range("N99").Value = -18958.522
# Row 135
$ws.Range("H135").Value = 992.5
$ws.Range("I135").Value = 788.44446
$ws.Range("J135").Value = 1359.8
$ws.Range("K135").Value = 7096.00014
$ws.Range("L135").Value = 12238.2
$ws.Range("M135").Value = -4561.00014
$ws.Range("N135").Value = -17308.2

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 9664.444
$ws.Range("I80").Value = 8930
$ws.Range("J80").Value = 11133.333
$ws.Range("K80").Value = 8930
$ws.Range("L80").Value = 11133.333
$ws.Range("M80").Value = -7932
$ws.Range("N80").Value = -13129.333
# Row 83
$ws.Range("H83").Value = 9664.444
$ws.Range("I83").Value = 8930
$ws.Range("J83").Value = 11133.333
$ws.Range("K83").Value = 44650
$ws.Range("L83").Value = 55666.665
$ws.Range("M83").Value = -39658
$ws.Range("N83").Value = -65650.66500000001
# Row 102
$ws.Range("H102").Value = 1322.1666
$ws.Range("I102").Value = 1322.1666
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1322.1666
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 299.8334
$ws.Range("N102").Value = ""

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 0
$ws.Range("I55").Value = 0
$ws.Range("K55").Value = 0
$ws.Range("M55").Value = ""
# Row 61
$ws.Range("H61").Value = 12752150
$ws.Range("I61").Value = 8502283
$ws.Range("J61").Value = 25501750
$ws.Range("K61").Value = 8502283
$ws.Range("L61").Value = 25501750
$ws.Range("M61").Value = -8502081
$ws.Range("N61").Value = -25502154
# Row 82
$ws.Range("H82").Value = 5259.778
$ws.Range("I82").Value = 1939.6666
$ws.Range("K82").Value = 1939.6666
$ws.Range("M82").Value = -1578.6666
# Row 85
$ws.Range("H85").Value = 5259.778
$ws.Range("I85").Value = 1939.6666
$ws.Range("K85").Value = 1939.6666
$ws.Range("M85").Value = -691.6666
# Row 93
$ws.Range("H93").Value = 0
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = ""
$ws.Range("N93").Value = ""
# Row 113
$ws.Range("H113").Value = 12752150
$ws.Range("I113").Value = 8502283
$ws.Range("J113").Value = 25501750
$ws.Range("K113").Value = 8502283
$ws.Range("L113").Value = 25501750
$ws.Range("M113").Value = -8500113
$ws.Range("N113").Value = -25506090

$ws = $wb.Worksheets.Item("WVR")
# Row 6
$ws.Range("H6").Value = 500
$ws.Range("I6").Value = 500
$ws.Range("K6").Value = 500
$ws.Range("M6").Value = -385
# Row 96
$ws.Range("H96").Value = 3372.5
$ws.Range("I96").Value = 3495
$ws.Range("K96").Value = 3495
$ws.Range("M96").Value = -2122
# Row 100
$ws.Range("H100").Value = 1789.6666
$ws.Range("I100").Value = 1582.6666
$ws.Range("K100").Value = 3165.3332
$ws.Range("M100").Value = -2624.3332
# Row 107
$ws.Range("H107").Value = 999.7143
$ws.Range("I107").Value = 999.6667
$ws.Range("K107").Value = 2999.0001
$ws.Range("M107").Value = -1079.0001
# Row 123
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = ""
# Row 126
$ws.Range("H126").Value = 2310.5454
$ws.Range("I126").Value = 2205.158
$ws.Range("J126").Value = 2978
$ws.Range("K126").Value = 6615.474
$ws.Range("L126").Value = 8934
$ws.Range("M126").Value = -4145.474
